# Generate Report for Archive
#
# 1. The shared string "Ready for handoff" becomes "In Translation".
#    It is referenced by every "Status" cell in the workbook (the
#    Overview sheet's zh-cn/de-de status columns, plus the Status column
#    on each per-locale detail sheet), so every one of those cells needs
#    to be rewritten to the new text.
# 2. The "Status" column is narrowed (from the wide ~17.22-char OOXML
#    width down to ~13.41 chars) on the Overview sheet (columns E and F)
#    and on the zh-cn / de-de sheets (column C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1. Status text: "Ready for handoff" -> "In Translation" ---------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- 2. Narrow the Status columns -------------------------------------------
# COM ColumnWidth of 12.5 round-trips through this host's pixel grid to the
# OOXML column width closest to the target 13.4101845877511.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
